# Update "Förändrad" date column (C) for rows 2-11 from 2023-09-16 (45185)
# to 2023-10-05 (45204), as part of an automatic data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..11) {
    $ws.Cells.Item($row, 3).Value = 45204
}
